$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 227, shifting rows 227-338 down to 228-339.
$ws.Rows.Item(227).Insert()

# Populate the newly inserted row 227 with the new weekly price observation.
$ws.Cells.Item(227, 1).Value = 10
$ws.Cells.Item(227, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(227, 3).Value = "La Araucanía"
$ws.Cells.Item(227, 4).Value = 45202
$ws.Cells.Item(227, 5).Value = 9
$ws.Cells.Item(227, 6).Value = 100112005
$ws.Cells.Item(227, 7).Value = "Puerro"
$ws.Cells.Item(227, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(227, 9).Value = "Primera"
$ws.Cells.Item(227, 10).Value = 65
$ws.Cells.Item(227, 11).Value = 8000
$ws.Cells.Item(227, 12).Value = 8000
$ws.Cells.Item(227, 13).Value = 8000
$ws.Cells.Item(227, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(227, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(227, 16).Value = 667
$ws.Cells.Item(227, 17).Value = 12
$ws.Cells.Item(227, 18).Value = "Hortaliza"
